$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "peliculas o documentales" - add two new rows (8 and 9)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("peliculas o documentales")

# Row 8: Sumergidos / Arriendo / Accion / HD / 2019 / $3.490
$ws1.Range("A8").Value = "Sumergidos"
$ws1.Range("C8").Value = "Arriendo"
$ws1.Range("D8").Value = "Acción"
$ws1.Range("E8").Value = "HD"
$ws1.Range("G8").Value = "$3.490"

# F8 must hold a genuine number (2019) while keeping the column's text
# number-format (style 1). Typing a number straight into a Text-formatted
# cell makes Excel store it as text, so reset to the default style first,
# enter the number, then paint the original text-format style back on.
$ws1.Range("F8").Style = "Normal"
$ws1.Range("F8").Value = 2019
$ws1.Range("F2").Copy()
$ws1.Range("F8").PasteSpecial(-4122)

# Row 9: This Is Us / FOX PREMIUM / Premium / Drama / HD / 2016
$ws1.Range("A9").Value = "This Is Us"
$ws1.Range("B9").Value = "FOX PREMIUM"
$ws1.Range("C9").Value = "Premium"
$ws1.Range("D9").Value = "Drama"
$ws1.Range("E9").Value = "HD"

$ws1.Range("F9").Style = "Normal"
$ws1.Range("F9").Value = 2016
$ws1.Range("F3").Copy()
$ws1.Range("F9").PasteSpecial(-4122)

# ---------------------------------------------------------------
# Sheet 2: "series" - fix row 4 year/episodes to real numbers and
# add a new row 5 for "My Brilliant Friend" Temporada 02
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("series")

# F4 ("2019" text) -> 2019 number
$ws2.Range("F4").Style = "Normal"
$ws2.Range("F4").Value = 2019
$ws2.Range("F2").Copy()
$ws2.Range("F4").PasteSpecial(-4122)

# H4 ("5" text) -> 5 number
$ws2.Range("H4").Style = "Normal"
$ws2.Range("H4").Value = 5
$ws2.Range("H2").Copy()
$ws2.Range("H4").PasteSpecial(-4122)

# New row 5
$ws2.Range("A5").Value = "My Brilliant Friend"
$ws2.Range("C5").Value = "Gratis"
$ws2.Range("D5").Value = "Drama"
$ws2.Range("E5").Value = "HD"
$ws2.Range("F5").Value = "2020"
$ws2.Range("G5").Value = "Temporada 02"
$ws2.Range("H5").Value = "1"
